# "Added first page index.html" -- jalon 1/Gestion du temps.xlsx
#
# 1) Feuil1 -> "Jalon 1"
# 2) new sheet "Jalon 2" appended after "Jalon 1", becomes the active sheet
# 3) "Jalon 2" gets a small time-tracking table (same layout/style as "Jalon 1")
# 4) selection/active-tab bookkeeping on both sheets

$wb = $excel.ActiveWorkbook

# --- sheet 1: rename -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Jalon 1"

# --- sheet 2: create, placed right after "Jalon 1" --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Jalon 2"

# --- header row (row 1) ------------------------------------------------
# Written so the brand-new shared strings get created in the same order
# they appear in the target workbook: B2, B3, C2, D1, C3, B4.
$ws2.Range("A1").Value = "demi-journée"
$ws2.Range("B1").Value = "description des tâches effectuées"
$ws2.Range("C1").Value = "Tâches unitaires"
$ws2.Range("E1").Value = "statut"
$ws2.Range("F1").Value = "durée idéale"
$ws2.Range("G1").Value = "Delta"

# --- data rows ----------------------------------------------------------
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "prise de connaissance du fil rouge"
$ws2.Range("B3").Value = "Rédaction de la maquette, ranger dossier et MAJ Git"
$ws2.Range("C2").Value = "Maquette"
$ws2.Range("D1").Value = "temps passé (en heures)"
$ws2.Range("C3").Value = "html"
$ws2.Range("B4").Value = "Rédaction de la maquette, début html"

$ws2.Range("D2").Value = 8
$ws2.Range("E2").Value = "WIP"
$ws2.Range("F2").Value = 8

$ws2.Range("A3").Value = 2
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = "WIP"
$ws2.Range("F3").Value = 8

$ws2.Range("A4").Value = 3

# --- formulas: Delta column ---------------------------------------------
$ws2.Range("G2").Formula = "=F2-D2"
$ws2.Range("G3:G12").Formula = "=F3-D3"

# --- header formatting: reuse "Jalon 1" header style (A1) --------------
$ws1.Range("A1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false
$ws2.Rows.Item(1).RowHeight = 27.6

# --- column B width (matches "Jalon 1" wrapping column look) -----------
$ws2.Columns.Item(2).ColumnWidth = 42.666666666666664

# --- selection bookkeeping -----------------------------------------------
$ws1.Range("A1:G1").Select()
$ws2.Activate()
$ws2.Range("D4").Select()
